$wb = $excel.ActiveWorkbook

# --- CAPEX sheet ---
$capex = $wb.Worksheets.Item("CAPEX")

# ADSL row (row 2): OPEX/total literal values reduced
$capex.Range("W2").Value = 10000
$capex.Range("X2").Value = 10000

# FTTC_GPON_25 (row 3): lower FIT (S column)
$capex.Range("S3").Value = 122492.30468757232

# FTTH_UDWDM_100 (row 6): inflated OPEX-related figures
$capex.Range("R6").Value = 758286.15266324603
$capex.Range("U6").Value = 160450
$capex.Range("V6").Value = 210910

# FTTH_XGPON_100 (row 7): R7 becomes literal text (not a formula)
$capex.Range("R7").Value = "384090.367674523+20*5000"

# FTTC_GPON_100 (row 8): lower FIT (S column)
$capex.Range("S8").Value = 122492.30468757232

# FTTC_Hybridpon_25 (row 11): lower FIT (S column)
$capex.Range("S11").Value = 81959.434330663411

# FTTC_Hybridpon_100 (row 13): R13 formula replaced by literal; FIT/CAPEX reduced
$capex.Range("R13").Value = 1012244.723172249
$capex.Range("T13").Value = 99287.142993292597
$capex.Range("V13").Value = 90438.2

# FTTB_Hybridpon_100 (row 15): R15 formula replaced by literal; FIT reduced
$capex.Range("R15").Value = 368464.72317224898
$capex.Range("T15").Value = 91410.786993292611

# New analysis block (rows 32-49): M = M+N+O of the row 29 rows earlier,
# N = P+Q+R of the row 29 rows earlier (N only through row 44)
for ($r = 32; $r -le 49; $r++) {
    $src = $r - 29
    $capex.Cells.Item($r, 13).Formula = "=M$src+N$src+O$src"
    if ($r -le 44) {
        $capex.Cells.Item($r, 14).Formula = "=P$src+Q$src+R$src"
    }
}

# --- OPEX sheet ---
$opex = $wb.Worksheets.Item("OPEX")

# ADSL (row 2): average formula replaced by a literal
$opex.Range("B2").Value = 10000

# FTTH_UDWDM_100 (row 6): inflated opex
$opex.Range("B6").Value = 19212.762708578561
